$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in column D whose new value looks like a plain number (e.g. "1.00",
# "0.0843") must be forced to text format first, otherwise Excel would
# auto-convert them to a numeric value and drop formatting like trailing
# zeros or leading zeros (the source data always stores these as text).
$textCells = @("D5","D6","D7","D10","D11","D14","D16","D19","D20","D22","D23","D24","D25","D29","D30","D33","D38","D41","D44","D45","D46","D47","D48","D49","D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Row 2 - Bitcoin
$ws.Range("D2").Value = "26.821.64"
$ws.Range("E2").Value = "  -0.10%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "1.638.55"
$ws.Range("E3").Value = "  -0.53%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  -0.82%  "

# Row 5 - BNB
$ws.Range("D5").Value = "218.86"
$ws.Range("E5").Value = "  +0.56%  "

# Row 6 - XRP
$ws.Range("D6").Value = "0.499"
$ws.Range("E6").Value = "  -0.97%  "

# Row 7 - USDC
$ws.Range("D7").Value = "1.00"
$ws.Range("E7").Value = "  -0.85%  "

# Row 8 - Cardano
$ws.Range("E8").Value = "  -0.51%  "

# Row 9 - Dogecoin
$ws.Range("E9").Value = "  -0.95%  "

# Row 10 - Solana
$ws.Range("D10").Value = "19.25"
$ws.Range("E10").Value = "  +0.24%  "

# Row 11 - TRON
$ws.Range("D11").Value = "0.0843"
$ws.Range("E11").Value = "  +0.02%  "

# Row 12 - WrappedliquidstakedEther2.0
$ws.Range("D12").Value = "1.865.23"
$ws.Range("E12").Value = "  -0.54%  "

# Row 13 - WrappedEther
$ws.Range("D13").Value = "1.639.83"
$ws.Range("E13").Value = "  -0.78%  "

# Row 14 - Polkadot
$ws.Range("D14").Value = "4.13"
$ws.Range("E14").Value = "  -1.33%  "

# Row 15 - Polygon
$ws.Range("E15").Value = "  -0.58%  "

# Row 16 - Litecoin
$ws.Range("D16").Value = "64.59"
$ws.Range("E16").Value = "  -0.27%  "

# Row 17 - WrappedBTC
$ws.Range("D17").Value = "26.809.77"
$ws.Range("E17").Value = "  -0.16%  "

# Row 18 - ShibaInu
$ws.Range("E18").Value = "  -0.63%  "

# Row 19 - BitcoinCash
$ws.Range("D19").Value = "215.01"
$ws.Range("E19").Value = "  +0.47%  "

# Row 20 - Dai
$ws.Range("D20").Value = "1.00"
$ws.Range("E20").Value = "  -0.81%  "

# Row 21 - Uniswap
$ws.Range("E21").Value = "  -0.05%  "

# Row 22 - Chainlink
$ws.Range("D22").Value = "6.28"
$ws.Range("E22").Value = "  -0.16%  "

# Row 23 - Toncoin
$ws.Range("D23").Value = "2.35"
$ws.Range("E23").Value = "  -2.44%  "

# Row 24 - Avalanche
$ws.Range("D24").Value = "9.09"
$ws.Range("E24").Value = "  -2.85%  "

# Row 25 - Monero
$ws.Range("D25").Value = "147.62"
$ws.Range("E25").Value = "  +1.62%  "

# Row 26 - BinanceUSD
$ws.Range("E26").Value = "  -0.74%  "

# Row 27 - Stellar
$ws.Range("E27").Value = "  +0.18%  "

# Row 28 - Cosmos
$ws.Range("E28").Value = "  -0.48%  "

# Row 29 - EthereumClassic
$ws.Range("D29").Value = "15.69"
$ws.Range("E29").Value = "  -0.04%  "

# Row 30 - Hedera
$ws.Range("D30").Value = "0.0504"
$ws.Range("E30").Value = "  -1.88%  "

# Row 31 - PancakeSwap
$ws.Range("E31").Value = "  +1.27%  "

# Row 32 - Filecoin
$ws.Range("E32").Value = "  +2.21%  "

# Row 33 - InternetComputer(DFINITY)
$ws.Range("D33").Value = "2.98"
$ws.Range("E33").Value = "  -0.15%  "

# Row 34 - LidoDAOToken
$ws.Range("E34").Value = "  +0.06%  "

# Row 35 - Maker
$ws.Range("D35").Value = "1.261.10"
$ws.Range("E35").Value = "  -1.44%  "

# Row 36 - HuobiToken
$ws.Range("E36").Value = "  -0.44%  "

# Row 37 - VeChain
$ws.Range("E37").Value = "  -0.13%  "

# Row 38 - ImmutableX
$ws.Range("D38").Value = "0.528"
$ws.Range("E38").Value = "  -1.61%  "

# Row 39 - ARBITRUM
$ws.Range("E39").Value = "  -1.70%  "

# Row 40 - PaxDollar
$ws.Range("E40").Value = "  -0.76%  "

# Row 41 - TrustWalletToken
$ws.Range("D41").Value = "0.807"
$ws.Range("E41").Value = "  -0.72%  "

# Row 42 - FraxShare
$ws.Range("E42").Value = "  -0.54%  "

# Row 43 - RocketPoolETH
$ws.Range("D43").Value = "1.775.13"
$ws.Range("E43").Value = "  -1.19%  "

# Row 44 - MXToken
$ws.Range("D44").Value = "2.13"
$ws.Range("E44").Value = "  -4.46%  "

# Row 45 - Quant
$ws.Range("D45").Value = "92.06"
$ws.Range("E45").Value = "  +0.53%  "

# Row 46 - Aave
$ws.Range("D46").Value = "60.05"
$ws.Range("E46").Value = "  +1.68%  "

# Row 47 - RenderToken
$ws.Range("D47").Value = "1.59"
$ws.Range("E47").Value = "  -0.52%  "

# Row 48 - was BabyDogeCoin, now Cronos
$ws.Range("B48").Value = "Cronos"
$ws.Range("C48").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D48").Value = "0.0516"
$ws.Range("E48").Value = "  -0.63%  "

# Row 49 - was Cronos, now EnergySwap
$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").Value = "7.55"
$ws.Range("E49").Value = "  -1.58%  "

# Row 50 - Algorand
$ws.Range("E50").Value = "  -1.40%  "

# Row 51 - USDD
$ws.Range("D51").Value = "1.00"
$ws.Range("E51").Value = "  -0.92%  "
